$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D9 text (shared string index 11 content change)
$ws.Range("D9").Value = "Apresenta lista de Serviços efectuados"

# Row 10: add flow step 3 - "Escolhe Serviço"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "Escolhe Serviço"

# Row 11: add flow step 4 - "Apresenta Avaliação"
$ws.Range("B11").Value = 4
$ws.Range("D11").Value = "Apresenta Avaliação"

# Update the active selection to D12
$ws.Range("D12").Select()
